$d = $word.ActiveDocument

# --- Apply the built-in "Heading 1" style to the (sole, empty) paragraph
#     in the document body. Word materialises the built-in "Heading1" /
#     "Heading1Char" style pair into styles.xml the first time the style
#     is actually used. ---
$p = $d.Paragraphs(1)
$p.Style = "Heading 1"

# Grab the freshly-minted paragraph style and square up the bits that
# the default built-in stub leaves out (name casing, next style, link
# to its companion run-level style, and the complex-script run size).
$heading1 = $d.Styles("Heading 1")
$heading1.NameLocal = "heading 1"
$heading1.NextParagraphStyle = "Normal"

# Create (and wire up) the linked character style "Heading 1 Char" that
# Word keeps alongside every heading paragraph style.
$heading1Char = $d.Styles.Add("Heading 1 Char", 2)
$heading1.LinkStyle = "Heading1Char"
$heading1Char.LinkStyle = "Heading1"
$heading1Char.BaseStyle = "DefaultParagraphFont"
$heading1Char.Priority = 9
$heading1Char.Font.Bold = $true
$heading1Char.Font.Size = 16
$heading1Char.Font.SizeBi = 16

# Complex-script size on the paragraph style itself (w:szCs).
$heading1.Font.SizeBi = 16
